$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 2448.6
$ws.Cells.Item(2, 9).Value = 199
$ws.Cells.Item(2, 11).Value = 199
$ws.Cells.Item(2, 13).Value = -86
# Row 74
$ws.Cells.Item(74, 8).Value = 8320
$ws.Cells.Item(74, 9).Value = 7800
$ws.Cells.Item(74, 11).Value = 7800
$ws.Cells.Item(74, 13).Value = -6864
# Row 77
$ws.Cells.Item(77, 8).Value = 8320
$ws.Cells.Item(77, 9).Value = 7800
$ws.Cells.Item(77, 11).Value = 39000
$ws.Cells.Item(77, 13).Value = -34320

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2782454.5
$ws.Cells.Item(32, 9).Value = 2917.4375
$ws.Cells.Item(32, 11).Value = 2917.4375
$ws.Cells.Item(32, 13).Value = -2630.4375
# Row 45
$ws.Cells.Item(45, 8).Value = 3555.7144
$ws.Cells.Item(45, 9).Value = 1972.5
$ws.Cells.Item(45, 11).Value = 1972.5
$ws.Cells.Item(45, 13).Value = -1595.5
# Row 61
$ws.Cells.Item(61, 8).Value = 2999.5
$ws.Cells.Item(61, 9).Value = 2999.5
$ws.Cells.Item(61, 11).Value = 2999.5
$ws.Cells.Item(61, 13).Value = -2787.5
# Row 136
$ws.Cells.Item(136, 8).Value = 2999.5
$ws.Cells.Item(136, 9).Value = 2999.5
$ws.Cells.Item(136, 11).Value = 8998.5
$ws.Cells.Item(136, 13).Value = -6448.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 1868.5
$ws.Cells.Item(99, 9).Value = 1958
$ws.Cells.Item(99, 11).Value = 1958
$ws.Cells.Item(99, 13).Value = -460

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5599.0713
$ws.Cells.Item(31, 9).Value = 3247.8667
$ws.Cells.Item(31, 11).Value = 3247.8667
$ws.Cells.Item(31, 13).Value = -2952.8667
# Row 34
$ws.Cells.Item(34, 8).Value = 5599.0713
$ws.Cells.Item(34, 9).Value = 3247.8667
$ws.Cells.Item(34, 11).Value = 3247.8667
$ws.Cells.Item(34, 13).Value = -3045.8667
# Row 75
$ws.Cells.Item(75, 8).Value = 39759.5
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 39759.5
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 39759.5
$ws.Cells.Item(75, 13).ClearContents()
$ws.Cells.Item(75, 14).Value = -41755.5
# Row 78
$ws.Cells.Item(78, 8).Value = 39759.5
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 39759.5
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 119278.5
$ws.Cells.Item(78, 13).ClearContents()
$ws.Cells.Item(78, 14).Value = -129262.5
# Row 99
$ws.Cells.Item(99, 8).Value = 2445.3333
$ws.Cells.Item(99, 9).Value = 1858.4286
$ws.Cells.Item(99, 10).Value = 4499.5
$ws.Cells.Item(99, 11).Value = 1858.4286
$ws.Cells.Item(99, 12).Value = 4499.5
$ws.Cells.Item(99, 13).Value = -360.4286
$ws.Cells.Item(99, 14).Value = -7495.5
# Row 105
$ws.Cells.Item(105, 8).Value = 1484.8
$ws.Cells.Item(105, 9).Value = 2262.5
$ws.Cells.Item(105, 10).Value = 966.3333
$ws.Cells.Item(105, 11).Value = 2262.5
$ws.Cells.Item(105, 12).Value = 966.3333
$ws.Cells.Item(105, 13).Value = -515.5
$ws.Cells.Item(105, 14).Value = -4460.3333
# Row 126
$ws.Cells.Item(126, 8).Value = 2445.3333
$ws.Cells.Item(126, 9).Value = 1858.4286
$ws.Cells.Item(126, 10).Value = 4499.5
$ws.Cells.Item(126, 11).Value = 5575.2858
$ws.Cells.Item(126, 12).Value = 13498.5
$ws.Cells.Item(126, 13).Value = -3105.2858
$ws.Cells.Item(126, 14).Value = -18438.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Cells.Item(9, 8).Value = 333.57144
$ws.Cells.Item(9, 10).Value = 616
$ws.Cells.Item(9, 12).Value = 1848
$ws.Cells.Item(9, 14).Value = -2296
# Row 46
$ws.Cells.Item(46, 8).Value = 904
$ws.Cells.Item(46, 9).Value = 904
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 2712
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -2621
$ws.Cells.Item(46, 14).ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 7222
$ws.Cells.Item(80, 9).Value = 10000
$ws.Cells.Item(80, 10).Value = 5833
$ws.Cells.Item(80, 11).Value = 10000
$ws.Cells.Item(80, 12).Value = 5833
$ws.Cells.Item(80, 13).Value = -9002
$ws.Cells.Item(80, 14).Value = -7829
# Row 83
$ws.Cells.Item(83, 8).Value = 7222
$ws.Cells.Item(83, 9).Value = 10000
$ws.Cells.Item(83, 10).Value = 5833
$ws.Cells.Item(83, 11).Value = 50000
$ws.Cells.Item(83, 12).Value = 29165
$ws.Cells.Item(83, 13).Value = -45008
$ws.Cells.Item(83, 14).Value = -39149
# Row 126
$ws.Cells.Item(126, 8).Value = 3629.5334
$ws.Cells.Item(126, 9).Value = 2632.111
$ws.Cells.Item(126, 11).Value = 7896.333
$ws.Cells.Item(126, 13).Value = -5426.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 9122.5
$ws.Cells.Item(7, 9).Value = 8499
$ws.Cells.Item(7, 10).Value = 9330.333000000001
$ws.Cells.Item(7, 11).Value = 8499
$ws.Cells.Item(7, 12).Value = 9330.333000000001
$ws.Cells.Item(7, 13).Value = -8387
$ws.Cells.Item(7, 14).Value = -9554.333000000001
# Row 16
$ws.Cells.Item(16, 8).Value = 185.77777
$ws.Cells.Item(16, 9).Value = 185.77777
$ws.Cells.Item(16, 11).Value = 185.77777
$ws.Cells.Item(16, 13).Value = -15.77777
# Row 26
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
# Row 93
$ws.Cells.Item(93, 8).Value = 1897
$ws.Cells.Item(93, 9).Value = 1862.6666
$ws.Cells.Item(93, 10).Value = 2000
$ws.Cells.Item(93, 11).Value = 1862.6666
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = -614.6666
$ws.Cells.Item(93, 14).Value = -4496
# Row 122
$ws.Cells.Item(122, 8).Value = 2868.6
$ws.Cells.Item(122, 9).Value = 2868.6
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 8605.799999999999
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -6155.799999999999
$ws.Cells.Item(122, 14).ClearContents()
# Row 126
$ws.Cells.Item(126, 8).Value = 9122.5
$ws.Cells.Item(126, 9).Value = 8499
$ws.Cells.Item(126, 10).Value = 9330.333000000001
$ws.Cells.Item(126, 11).Value = 25497
$ws.Cells.Item(126, 12).Value = 27990.999
$ws.Cells.Item(126, 13).Value = -23027
$ws.Cells.Item(126, 14).Value = -32930.999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Cells.Item(15, 8).Value = 10001
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
# Row 62
$ws.Cells.Item(62, 8).Value = 7870.3335
$ws.Cells.Item(62, 10).Value = 10600
$ws.Cells.Item(62, 12).Value = 10600
$ws.Cells.Item(62, 14).Value = -11848
# Row 65
$ws.Cells.Item(65, 8).Value = 7870.3335
$ws.Cells.Item(65, 10).Value = 10600
$ws.Cells.Item(65, 12).Value = 53000
$ws.Cells.Item(65, 14).Value = -59240
# Row 93
$ws.Cells.Item(93, 8).Value = 31694.5
$ws.Cells.Item(93, 10).Value = 31694.5
$ws.Cells.Item(93, 12).Value = 31694.5
$ws.Cells.Item(93, 14).Value = -36686.5
# Row 118
$ws.Cells.Item(118, 8).Value = 55000
$ws.Cells.Item(118, 10).Value = 55000
$ws.Cells.Item(118, 12).Value = 55000
$ws.Cells.Item(118, 14).Value = -58314
